$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 133 ("photograph" / "写真|しゃしん") was removed from the wordlist.
# Deleting the entire row shifts rows 134:146 up to 133:145,
# matching the diff (dimension shrinks from B146 to B145).
$ws.Rows.Item(133).Delete()
